# Apply updates to Fruta / hortaliza, semanal (Granada rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44664
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 21600
$ws.Range("O2").Value = 21600
$ws.Range("P2").Value = 21600
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1200
# Row 3
$ws.Range("D3").Value = 44664
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1000
# Row 4
$ws.Range("D4").Value = 44664
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 889
# Row 5
$ws.Range("D5").Value = 44698
$ws.Range("K5").Value = "Wonderfull"
$ws.Range("M5").Value = 280
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("S5").Value = 1200
# Row 6
$ws.Range("D6").Value = 44698
$ws.Range("K6").Value = "Wonderfull"
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 1000
# Row 7
$ws.Range("D7").Value = 44698
$ws.Range("K7").Value = "Wonderfull"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 800
# Row 8
$ws.Range("D8").Value = 44678
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 290
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1000
# Row 9
$ws.Range("D9").Value = 44678
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 220
$ws.Range("R9").Value = "Región de O'Higgins"
# Row 10
$ws.Range("D10").Value = 44309
$ws.Range("K10").Value = "Wonderfull"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 40
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("Q10").Value = "$/caja 15 kilos granel"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("T10").Value = 15
# Row 11
$ws.Range("D11").Value = 44309
$ws.Range("K11").Value = "Wonderfull"
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 70
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("Q11").Value = "$/caja 15 kilos granel"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 15
# Row 12
$ws.Range("D12").Value = 44285
$ws.Range("K12").Value = "Wonderfull"
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("Q12").Value = "$/caja 15 kilos empedrada"
$ws.Range("R12").Value = "Provincia del Elquí"
$ws.Range("S12").Value = 1200
$ws.Range("T12").Value = 15
# Row 13
$ws.Range("D13").Value = 44285
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 90
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("Q13").Value = "$/caja 15 kilos empedrada"
$ws.Range("R13").Value = "Provincia del Elquí"
$ws.Range("S13").Value = 1000
$ws.Range("T13").Value = 15
# Row 14
$ws.Range("D14").Value = 44285
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 75
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = "$/caja 15 kilos empedrada"
$ws.Range("R14").Value = "Provincia del Elquí"
$ws.Range("S14").Value = 800
$ws.Range("T14").Value = 15
# Row 15
$ws.Range("D15").Value = 44644
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("M15").Value = 180
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 18000
$ws.Range("Q15").Value = "$/caja 15 kilos granel"
$ws.Range("R15").Value = "Provincia de Limarí"
$ws.Range("T15").Value = 15
# Row 16
$ws.Range("D16").Value = 44644
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("N16").Value = 13500
$ws.Range("O16").Value = 13500
$ws.Range("P16").Value = 13500
$ws.Range("Q16").Value = "$/caja 15 kilos granel"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 900
$ws.Range("T16").Value = 15
# Row 17
$ws.Range("D17").Value = 44644
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("M17").Value = 290
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = "$/caja 15 kilos granel"
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("T17").Value = 15
# Row 20
$ws.Range("D20").Value = 44687
$ws.Range("K20").Value = "Wonderfull"
$ws.Range("M20").Value = 220
$ws.Range("N20").Value = 21000
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 21000
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 1167
# Row 21
$ws.Range("D21").Value = 44687
$ws.Range("K21").Value = "Wonderfull"
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 15000
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 833
# Row 22
$ws.Range("D22").Value = 44687
$ws.Range("K22").Value = "Wonderfull"
$ws.Range("M22").Value = 280
$ws.Range("N22").Value = 10000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 10000
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 556
# Row 23
$ws.Range("D23").Value = 44694
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 21600
$ws.Range("O23").Value = 21600
$ws.Range("P23").Value = 21600
$ws.Range("Q23").Value = "$/caja 18 kilos granel"
$ws.Range("R23").Value = "Región de O'Higgins"
$ws.Range("T23").Value = 18
# Row 24
$ws.Range("D24").Value = 44694
$ws.Range("M24").Value = 220
$ws.Range("N24").Value = 18000
$ws.Range("O24").Value = 18000
$ws.Range("P24").Value = 18000
$ws.Range("Q24").Value = "$/caja 18 kilos granel"
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("T24").Value = 18
# Row 25
$ws.Range("D25").Value = 44694
$ws.Range("M25").Value = 250
$ws.Range("N25").Value = 14400
$ws.Range("O25").Value = 14400
$ws.Range("P25").Value = 14400
$ws.Range("Q25").Value = "$/caja 18 kilos granel"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("T25").Value = 18
# Row 29
$ws.Range("D29").Value = 44706
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 16000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 16000
$ws.Range("Q29").Value = "$/caja 18 kilos granel"
$ws.Range("S29").Value = 889
$ws.Range("T29").Value = 18
# Row 30
$ws.Range("D30").Value = 44706
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 220
$ws.Range("N30").Value = 12500
$ws.Range("O30").Value = 12500
$ws.Range("P30").Value = 12500
$ws.Range("Q30").Value = "$/caja 18 kilos granel"
$ws.Range("S30").Value = 694
$ws.Range("T30").Value = 18
# Row 31
$ws.Range("D31").Value = 44649
$ws.Range("K31").Value = "Sin especificar"
$ws.Range("M31").Value = 220
$ws.Range("N31").Value = 21600
$ws.Range("O31").Value = 21600
$ws.Range("P31").Value = 21600
$ws.Range("Q31").Value = "$/caja 18 kilos granel"
$ws.Range("R31").Value = "Provincia de Limarí"
$ws.Range("T31").Value = 18
# Row 32
$ws.Range("D32").Value = 44649
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("M32").Value = 250
$ws.Range("N32").Value = 16200
$ws.Range("O32").Value = 16200
$ws.Range("P32").Value = 16200
$ws.Range("Q32").Value = "$/caja 18 kilos granel"
$ws.Range("R32").Value = "Provincia de Limarí"
$ws.Range("S32").Value = 900
$ws.Range("T32").Value = 18
# Row 33
$ws.Range("D33").Value = 44649
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("M33").Value = 180
$ws.Range("N33").Value = 14400
$ws.Range("O33").Value = 14400
$ws.Range("P33").Value = 14400
$ws.Range("Q33").Value = "$/caja 18 kilos granel"
$ws.Range("R33").Value = "Provincia de Limarí"
$ws.Range("T33").Value = 18
# Row 34
$ws.Range("D34").Value = 44658
$ws.Range("M34").Value = 280
# Row 35
$ws.Range("D35").Value = 44658
$ws.Range("M35").Value = 330
# Row 36
$ws.Range("D36").Value = 44658
